# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (B16:J28) are re-sorted from descending period
# order (2103 .. 2003) to ascending period order (2003 .. 2103), carrying
# each row's "Valor Mora" (column F) along with its matching period so the
# period -> valor mora relationship is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 16
$endRow = 28

# Read current "Periodo Mora" (E) and "Valor Mora" (F) pairs for the data rows.
$periods = @()
$valores = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $periods += $ws.Range("E$r").Value()
    $valores += $ws.Range("F$r").Value()
}

# Build period -> valor mora lookup, then sort periods ascending.
$map = @{}
for ($i = 0; $i -lt $periods.Length; $i++) {
    $map[[string]$periods[$i]] = $valores[$i]
}
$sortedPeriods = $periods | Sort-Object { [int]$_ }

# Write the periods (and their corresponding valor mora) back in ascending order.
for ($i = 0; $i -lt $sortedPeriods.Length; $i++) {
    $r = $startRow + $i
    $period = $sortedPeriods[$i]
    $ws.Range("E$r").Value = $period
    $ws.Range("F$r").Value = $map[[string]$period]
}
